$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values in rows 2, 4, 5, 9, 13 ---
$ws.Range("G2").Value = 2.15
$ws.Range("H2").Value = 3.1
$ws.Range("J2").Value = 3
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6
$ws.Range("Z2").Value = 19
$ws.Range("AC2").Value = 6.5
$ws.Range("AF2").Value = 67
$ws.Range("AL2").Value = 34
$ws.Range("AM2").Value = 41
$ws.Range("AQ2").Value = 41
$ws.Range("AS2").Value = 251
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 9
$ws.Range("AY2").Value = 34
$ws.Range("BB2").Value = 351
$ws.Range("G4").Value = 4
$ws.Range("I4").Value = 1.91
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.63
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.25
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.8
$ws.Range("X4").Value = 19
$ws.Range("Y4").Value = 13
$ws.Range("AA4").Value = 34
$ws.Range("AC4").Value = 9
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 351
$ws.Range("AI4").Value = 8.5
$ws.Range("AK4").Value = 17
$ws.Range("AM4").Value = 29
$ws.Range("AR4").Value = 101
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.63
$ws.Range("AV4").Value = 51
$ws.Range("AX4").Value = 11
$ws.Range("AZ4").Value = 41
$ws.Range("BB4").Value = 151
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67
$ws.Range("BD9").Value = 151
$ws.Range("G13").Value = 2.15
$ws.Range("I13").Value = 3.3
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 9
$ws.Range("AC13").Value = 9
$ws.Range("AG13").Value = 301
$ws.Range("AI13").Value = 17
$ws.Range("AL13").Value = 29
$ws.Range("AN13").Value = 4
$ws.Range("AZ13").Value = 67

# --- Insert a new match row at row 16 (Paraguay - Primera Division: Cerro Porteno vs Guarani) ---
# This shifts the existing rows 16-25 down to rows 17-26.
$ws.Rows(16).Insert()

$ws.Range("A16").Value = "Q7FCwn13"
$ws.Range("B16").Value = "16/11/2024"
$ws.Range("C16").Value = "20:30"
$ws.Range("D16").Value = "PARAGUAY - PRIMERA DIVISION"
$ws.Range("E16").Value = "Cerro Porteno"
$ws.Range("F16").Value = "Guarani"
$ws.Range("G16").Value = 1.53
$ws.Range("H16").Value = 3.75
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 2.2
$ws.Range("K16").Value = 2.1
$ws.Range("L16").Value = 7
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 9
$ws.Range("O16").Value = 1.36
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 2.2
$ws.Range("R16").Value = 1.65
$ws.Range("S16").Value = 1.5
$ws.Range("T16").Value = 2.5
$ws.Range("U16").Value = 2.25
$ws.Range("V16").Value = 1.57
$ws.Range("W16").Value = 5.5
$ws.Range("X16").Value = 6.5
$ws.Range("Y16").Value = 9
$ws.Range("Z16").Value = 10
$ws.Range("AA16").Value = 15
$ws.Range("AB16").Value = 34
$ws.Range("AC16").Value = 8
$ws.Range("AD16").Value = 7.5
$ws.Range("AE16").Value = 21
$ws.Range("AF16").Value = 81
$ws.Range("AG16").Value = 201
$ws.Range("AH16").Value = 13
$ws.Range("AI16").Value = 34
$ws.Range("AJ16").Value = 21
$ws.Range("AK16").Value = 81
$ws.Range("AL16").Value = 51
$ws.Range("AM16").Value = 51
$ws.Range("AN16").Value = 3.25
$ws.Range("AO16").Value = 8
$ws.Range("AP16").Value = 23
$ws.Range("AQ16").Value = 26
$ws.Range("AR16").Value = 51
$ws.Range("AS16").Value = 201
$ws.Range("AT16").Value = 2.5
$ws.Range("AU16").Value = 10
$ws.Range("AV16").Value = 81
$ws.Range("AW16").Value = 7.5
$ws.Range("AX16").Value = 41
$ws.Range("AY16").Value = 41
$ws.Range("AZ16").Value = 151
$ws.Range("BA16").Value = 201
$ws.Range("BB16").Value = 351
$ws.Range("BC16").Value = 51
$ws.Range("BD16").Value = 51

Write-Host "Edit complete"
